$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.735750333333333
$ws.Range("H2").Value = 5.207250999999999
$ws.Range("I2").Value = 0.01516614046792869
$ws.Range("J2").Value = 0.01516614046792869
$ws.Range("M2").Value = 0.4890553333333333
$ws.Range("N2").Value = 1.467166
$ws.Range("O2").Value = 0.9644476581758422
$ws.Range("P2").Value = 0.9644476581758422
$ws.Range("Q2").Value = 0.8488779578517777
$ws.Range("R2").Value = 7.639901620665999
$ws.Range("S2").Value = 0.0146269486578597
$ws.Range("T2").Value = 0.0146269486578597

# Row 3
$ws.Range("G3").Value = 1.735750333333333
$ws.Range("H3").Value = 5.207250999999999
$ws.Range("I3").Value = 0.01516614046792869
$ws.Range("J3").Value = 0.01516614046792869
$ws.Range("O3").Value = 0.03555234182415776
$ws.Range("P3").Value = 0.03555234182415776
$ws.Range("Q3").Value = 0.03129210700933333
$ws.Range("R3").Value = 0.281628963084
$ws.Range("S3").Value = 0.0005391918100689929
$ws.Range("T3").Value = 0.0005391918100689927

# Row 4
$ws.Range("I4").Value = 0.8822596377334645
$ws.Range("J4").Value = 0.8822596377334644
$ws.Range("M4").Value = 0.4890553333333333
$ws.Range("N4").Value = 1.467166
$ws.Range("O4").Value = 0.9644476581758422
$ws.Range("P4").Value = 0.9644476581758422
$ws.Range("Q4").Value = 49.38176335356845
$ws.Range("R4").Value = 444.435870182116
$ws.Range("S4").Value = 0.8508932415151067
$ws.Range("T4").Value = 0.8508932415151066

# Row 5
$ws.Range("I5").Value = 0.8822596377334645
$ws.Range("J5").Value = 0.8822596377334644
$ws.Range("O5").Value = 0.03555234182415776
$ws.Range("P5").Value = 0.03555234182415776
$ws.Range("S5").Value = 0.03136639621835773
$ws.Range("T5").Value = 0.03136639621835773

# Row 6
$ws.Range("I6").Value = 0.1025742217986069
$ws.Range("J6").Value = 0.1025742217986069
$ws.Range("M6").Value = 0.4890553333333333
$ws.Range("N6").Value = 1.467166
$ws.Range("O6").Value = 0.9644476581758422
$ws.Range("P6").Value = 0.9644476581758422
$ws.Range("Q6").Value = 5.741275844884
$ws.Range("R6").Value = 51.67148260395599
$ws.Range("S6").Value = 0.09892746800287584
$ws.Range("T6").Value = 0.09892746800287583

# Row 7
$ws.Range("I7").Value = 0.1025742217986069
$ws.Range("J7").Value = 0.1025742217986069
$ws.Range("O7").Value = 0.03555234182415776
$ws.Range("P7").Value = 0.03555234182415776
$ws.Range("S7").Value = 0.003646753795731047
$ws.Range("T7").Value = 0.003646753795731047
